# Correction of the reference values (fmod / R) in the roughness test data.
# Columns A (fmod) and B (R) are stored as text (shared-string) cells, not
# numbers, so each numeric-looking value is written back as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New fmod (col A) values for rows 3..28
$fmod = @(
    "14.026899", "16.03378",  "19.51608",  "24.1325",   "27.152908",
    "33.18282",  "42.67351",  "48.009727", "51.526123", "55.798893",
    "58.02824",  "63.2473",   "68.39307",  "73.18517",  "78.91965",
    "85.31508",  "94.640015", "97.87623",  "112.84719", "116.1008",
    "134.20557", "142.05911", "176.54396", "208.32822", "247.7417",
    "325.8978"
)

# New R (col B) values for rows 3..28
$rvals = @(
    "0.120577574", "0.14782797", "0.1944676",  "0.26501298", "0.31240755",
    "0.420769",    "0.5939494",  "0.68385667", "0.7424558",  "0.78007925",
    "0.7985556",   "0.8280814",  "0.8409192",  "0.8384773",  "0.8165356",
    "0.77257395",  "0.68996763", "0.65123534", "0.5195601",  "0.49169108",
    "0.3902205",   "0.3513164",  "0.2457958",  "0.18558311", "0.13721035",
    "0.08642737"
)

for ($i = 0; $i -lt $fmod.Length; $i++) {
    $row = $i + 3

    $cellA = $ws.Range("A" + $row)
    $cellA.NumberFormat = "@"
    $cellA.Value = $fmod[$i]
    $cellA.Style = "Normal"

    $cellB = $ws.Range("B" + $row)
    $cellB.NumberFormat = "@"
    $cellB.Value = $rvals[$i]
    $cellB.Style = "Normal"
}

# Keep the sheet's recorded selection in sync with the now-larger data range
# (matches the workbook's saved view state after the edit).
$ws.Range("B3:B28").Select()
